$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that already have a species filled in (these only need a casing fix
# on column F; everything else about them is left untouched).
$specialRows = @{ 9 = "Squirrel"; 20 = "Other"; 44 = "Mink" }

for ($r = 2; $r -le 46; $r++) {
    if ($specialRows.ContainsKey($r)) {
        $ws.Cells.Item($r, 6).Value = $specialRows[$r]
    } else {
        $ws.Cells.Item($r, 2).Value = "Yes"
        $ws.Cells.Item($r, 6).Value = "Na"
        $ws.Cells.Item($r, 10).Value = 1
    }
}
